$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item("Overview")

# Update cell values
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "zh-cn"
$ws.Range("C1").Value = "de-de"
$ws.Range("D1").Value = "Latest Handoff Date"
$ws.Range("A2").Value = "ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-23 09:58:28"
$ws.Range("A3").Value = "ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-23 09:58:28"
$ws.Range("A4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-23 10:02:46"

# Rebuild hyperlinks (same target URLs, refreshed display text matching new cell content)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md", [Type]::Missing, [Type]::Missing, "ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md", [Type]::Missing, [Type]::Missing, "ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md") | Out-Null

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "File Extension"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Latest Handoff File"
$ws.Range("E1").Value = "Latest Handoff Datetime"
$ws.Range("F1").Value = "Latest Target File"
$ws.Range("G1").Value = "Latest Handback File"
$ws.Range("H1").Value = "Latest Handback DateTime"
$ws.Range("J1").Value = "Handoff Reason"
$ws.Range("K1").Value = "Dependency From"
$ws.Range("L1").Value = "Error Detail"
$ws.Range("A2").Value = "ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-23 09:58:20"
$ws.Range("F2").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md"
$ws.Range("G2").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-23 09:58:57"
$ws.Range("J2").Value = "Include"
$ws.Range("A3").Value = "ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-23 09:58:20"
$ws.Range("F3").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md"
$ws.Range("G3").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-23 09:58:57"
$ws.Range("J3").Value = "Include"
$ws.Range("A4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-23 10:02:38"
$ws.Range("F4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md"
$ws.Range("G4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-23 10:01:36"
$ws.Range("J4").Value = "Include"

# Rebuild hyperlinks (same target URLs, refreshed display text matching new cell content)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md", [Type]::Missing, [Type]::Missing, "ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dbc3fc6a9f3713ff3dfc3e5cef4b72a57bbaed9c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/43ae5a2a75c6a3a2ff67ab26f83888fe8360666d/e2e/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1ca6dd2a6dd135526c2757c9bdc8a88db3ac3847/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md", [Type]::Missing, [Type]::Missing, "ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2431d9a335e0608258005caaea302e10c4ba5a6f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/7ea068a0c538a0f0994510b7d6fda5870bf46805/e2e/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6864d2b55031a900f6699d9622d1d81f0fecbd1e/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2431d9a335e0608258005caaea302e10c4ba5a6f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/7ea068a0c538a0f0994510b7d6fda5870bf46805/e2e/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6864d2b55031a900f6699d9622d1d81f0fecbd1e/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.zh-cn.xlf") | Out-Null

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item("de-de")

# Update cell values
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "File Extension"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Latest Handoff File"
$ws.Range("E1").Value = "Latest Handoff Datetime"
$ws.Range("F1").Value = "Latest Target File"
$ws.Range("G1").Value = "Latest Handback File"
$ws.Range("H1").Value = "Latest Handback DateTime"
$ws.Range("J1").Value = "Handoff Reason"
$ws.Range("K1").Value = "Dependency From"
$ws.Range("L1").Value = "Error Detail"
$ws.Range("A2").Value = "ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf"
$ws.Range("E2").Value = "2016-03-23 09:58:28"
$ws.Range("F2").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md"
$ws.Range("G2").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf"
$ws.Range("H2").Value = "2016-03-23 09:59:11"
$ws.Range("J2").Value = "Include"
$ws.Range("A3").Value = "ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf"
$ws.Range("E3").Value = "2016-03-23 09:58:28"
$ws.Range("F3").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md"
$ws.Range("G3").Value = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf"
$ws.Range("H3").Value = "2016-03-23 09:59:11"
$ws.Range("J3").Value = "Include"
$ws.Range("A4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.de-de.xlf"
$ws.Range("E4").Value = "2016-03-23 10:02:46"
$ws.Range("F4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md"
$ws.Range("G4").Value = "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.de-de.xlf"
$ws.Range("H4").Value = "2016-03-23 10:01:58"
$ws.Range("J4").Value = "Include"

# Rebuild hyperlinks (same target URLs, refreshed display text matching new cell content)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md", [Type]::Missing, [Type]::Missing, "ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97c5cb9ed275f00b2f004017d51798f120af6484/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.de-de.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/bb358acf650bf3688d30cfe16a00eef5784e71f5/e2e/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5b8ae424355e53625b4ccceaeaf08f982f87e254/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.de-de.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/ffffa7d645dc-4b9c-4f62-983b-481d95446e89.md", [Type]::Missing, [Type]::Missing, "ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/deb48e9cd88a7a33eeae8433bdeae550f29df94b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/14d2c70d0160eb6158ec111ccefd2d4a3f32995a/e2e/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b5dd8ee7589011fb4d625a0ac4bdbf25f3375ea/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf", [Type]::Missing, [Type]::Missing, "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f4446313cdae6366cb78bc899a7e9231c8bce02b/e2e/ffffff1d46648d-714e-46ef-bb1b-fe2b283323ce.md", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/deb48e9cd88a7a33eeae8433bdeae550f29df94b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/14d2c70d0160eb6158ec111ccefd2d4a3f32995a/e2e/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b5dd8ee7589011fb4d625a0ac4bdbf25f3375ea/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf", [Type]::Missing, [Type]::Missing, "f61812b9-60ff-45eb-a28b-e2d51f8e3ae6.215739689dd5c3503497a43adde29844227f268e.de-de.xlf") | Out-Null
